$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp shown in F1
$ws.Range("F1").Value = "Last status check on: 06.02.2022 02:30"

# Row 3 (Tesco) price refresh: new current price, previous price becomes
# the "old" price, delta is now recorded as text, and the timestamp is
# now recorded as a plain text string instead of a formatted serial date.
$ws.Range("B3").Value = 35.9
$ws.Range("C3").Value = 35.5

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "+0.4"
$dCell.Style = "Normal"

$eCell = $ws.Range("E3")
$eCell.Value = "2022-02-06 02:31:09"
$eCell.Style = "Normal"
